$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update N column (MACRO_SCORE) for rows 2-7
$ws.Range("N2").Value = 54.86376272656823
$ws.Range("N3").Value = 54.86376272656823
$ws.Range("N4").Value = 54.86376272656823
$ws.Range("N5").Value = 54.86376272656823
$ws.Range("N6").Value = 54.86376272656823
$ws.Range("N7").Value = 54.86376272656823

# Row 3 (SamsungElec) updated close/RSI/5d return
$ws.Range("D3").Value = 107100
$ws.Range("E3").Value = 60.1
$ws.Range("F3").Value = 6.57

# Row 5 (SK hynix) updated close/RSI/5d return
$ws.Range("D5").Value = 535000
$ws.Range("E5").Value = 31.9
$ws.Range("F5").Value = 0.9399999999999999

# Row 6 and Row 7 swap content (240810.KS moves to row 6, DB HiTek moves to row 7),
# with row 7 (DB HiTek) getting freshly updated numeric values.
$ws.Range("B6").Value = "240810.KS,0P00017YB3,330568"
$ws.Range("C6").Value = "240810.KS"
$ws.Range("D6").Value = 61000
$ws.Range("E6").Value = 36.7
$ws.Range("F6").Value = 1.16
$ws.Range("H6").Value = 60
$ws.Range("I6").Value = 46
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 40.9

$ws.Range("B7").Value = "DB HiTek"
$ws.Range("C7").Value = "000990.KS"
$ws.Range("D7").Value = 64100
$ws.Range("E7").Value = 31.8
$ws.Range("F7").Value = 0.79
$ws.Range("H7").Value = 46
$ws.Range("I7").Value = 43
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 39.7

$wb.Save()
